# DEV 14 - Final Fix
# - Correct the PROJECT ID value for the third enquiry row (C3: 1 -> 4)
# - Standardise the date/time number format used by the ENQUIRY_DATE /
#   REPLY_DATE columns on rows 3 and 4 so they match the format already
#   used on row 2 (collapses the redundant duplicate date numFmt/style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the PROJECT ID for row 3
$ws.Range("C3").Value = 4

# Re-apply the canonical date/time format (same one used by F2:G2) to the
# enquiry/reply date cells on rows 3 and 4, so they stop using the stray
# duplicate number format.
$ws.Range("F3:G3").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Range("F4").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"

# Leave the selection where the author left it when they saved the file.
$ws.Range("F16").Select()
